$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 24667.223
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 27563.125
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 27563.125
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -28215.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 29933.334
$ws.Range("J63").Value = 29933.334
$ws.Range("L63").Value = 29933.334
$ws.Range("N63").Value = -31181.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 29933.334
$ws.Range("J66").Value = 29933.334
$ws.Range("L66").Value = 89800.00199999999
$ws.Range("N66").Value = -96040.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 896.6667
$ws.Range("I45").Value = 896.6667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 896.6667
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -519.6667
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 29666.666
$ws.Range("J64").Value = 29666.666
$ws.Range("L64").Value = 29666.666
$ws.Range("N64").Value = -30162.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 29666.666
$ws.Range("J67").Value = 29666.666
$ws.Range("L67").Value = 29666.666
$ws.Range("N67").Value = -31382.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5562.4443
$ws.Range("I122").Value = 7430.6665
$ws.Range("J122").Value = 1826
$ws.Range("K122").Value = 22291.9995
$ws.Range("L122").Value = 5478
$ws.Range("M122").Value = -19841.9995
$ws.Range("N122").Value = -10378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 26418
$ws.Range("J123").Value = 26418
$ws.Range("L123").Value = 26418
$ws.Range("N123").Value = -36218

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 78943.62
$ws.Range("I86").Value = 2014
$ws.Range("J86").Value = 168694.83
$ws.Range("K86").Value = 2014
$ws.Range("L86").Value = 168694.83
$ws.Range("M86").Value = -891
$ws.Range("N86").Value = -170940.83

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 78943.62
$ws.Range("I89").Value = 2014
$ws.Range("J89").Value = 168694.83
$ws.Range("K89").Value = 10070
$ws.Range("L89").Value = 843474.1499999999
$ws.Range("M89").Value = -4454
$ws.Range("N89").Value = -854706.1499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3399.077
$ws.Range("I105").Value = 3386.9
$ws.Range("K105").Value = 3386.9
$ws.Range("M105").Value = -1639.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 369.95
$ws.Range("I22").Value = 217.58824
$ws.Range("J22").Value = 1233.3334
$ws.Range("K22").Value = 217.58824
$ws.Range("L22").Value = 1233.3334
$ws.Range("M22").Value = 132.41176
$ws.Range("N22").Value = -1933.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1753.9412
$ws.Range("I31").Value = 1143.3556
$ws.Range("J31").Value = 6333.3335
$ws.Range("K31").Value = 1143.3556
$ws.Range("L31").Value = 6333.3335
$ws.Range("M31").Value = -848.3556000000001
$ws.Range("N31").Value = -6923.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1753.9412
$ws.Range("I34").Value = 1143.3556
$ws.Range("J34").Value = 6333.3335
$ws.Range("K34").Value = 1143.3556
$ws.Range("L34").Value = 6333.3335
$ws.Range("M34").Value = -941.3556000000001
$ws.Range("N34").Value = -6737.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 44.95238
$ws.Range("I12").Value = 7.6
$ws.Range("J12").Value = 56.625
$ws.Range("K12").Value = 22.8
$ws.Range("L12").Value = 169.875
$ws.Range("M12").Value = 150.2
$ws.Range("N12").Value = -515.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5249.933
$ws.Range("I133").Value = 1770.2
$ws.Range("J133").Value = 6989.8
$ws.Range("K133").Value = 5310.6
$ws.Range("L133").Value = 20969.4
$ws.Range("M133").Value = -250.6000000000004
$ws.Range("N133").Value = -31089.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3509
$ws.Range("I134").Value = 2122
$ws.Range("J134").Value = 5242.75
$ws.Range("K134").Value = 6366
$ws.Range("L134").Value = 15728.25
$ws.Range("M134").Value = -1296
$ws.Range("N134").Value = -25868.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3874.0667
$ws.Range("I80").Value = 3733.889
$ws.Range("J80").Value = 4084.3333
$ws.Range("K80").Value = 3733.889
$ws.Range("L80").Value = 4084.3333
$ws.Range("M80").Value = -2735.889
$ws.Range("N80").Value = -6080.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3874.0667
$ws.Range("I83").Value = 3733.889
$ws.Range("J83").Value = 4084.3333
$ws.Range("K83").Value = 18669.445
$ws.Range("L83").Value = 20421.6665
$ws.Range("M83").Value = -13677.445
$ws.Range("N83").Value = -30405.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3839.2354
$ws.Range("I122").Value = 3387.5
$ws.Range("J122").Value = 4240.778
$ws.Range("K122").Value = 10162.5
$ws.Range("L122").Value = 12722.334
$ws.Range("M122").Value = -7712.5
$ws.Range("N122").Value = -17622.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 37648
$ws.Range("J131").Value = 37648
$ws.Range("L131").Value = 37648
$ws.Range("N131").Value = -47728

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1440.6666
$ws.Range("I16").Value = 1495.75
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1495.75
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1325.75
$ws.Range("N16").Value = -1340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 534
$ws.Range("I55").Value = 300
$ws.Range("K55").Value = 300
$ws.Range("M55").Value = -127

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22364.348
$ws.Range("J123").Value = 22364.348
$ws.Range("L123").Value = 22364.348
$ws.Range("N123").Value = -32164.348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 60048.75
$ws.Range("J125").Value = 60048.75
$ws.Range("L125").Value = 60048.75
$ws.Range("N125").Value = -69888.75
